# Shift the dates in column F (rows 2-7) forward by 4 days.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 44592
$ws.Range("F3").Value = 44591
$ws.Range("F4").Value = 44590
$ws.Range("F5").Value = 44589
$ws.Range("F6").Value = 44588
$ws.Range("F7").Value = 44587
